# Populate MPA transfer-scenario id columns (K, L, N, O) on the "Data" sheet
# with the real numeric scenario/step ids used by the test-automation upload,
# replacing the old placeholder scenario-name strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Cells.Item(6, 11).Value = 60000168
$ws.Cells.Item(6, 14).Value = ""

$ws.Cells.Item(7, 11).Value = 60000000
$ws.Cells.Item(7, 12).Value = 125
$ws.Cells.Item(7, 14).Value = 60000169

$ws.Cells.Item(8, 11).Value = 60000000
$ws.Cells.Item(8, 12).Value = 125
$ws.Cells.Item(8, 14).Value = 60000000
$ws.Cells.Item(8, 15).Value = 126

$ws.Cells.Item(9, 11).Value = 60000168
$ws.Cells.Item(9, 14).Value = ""

$ws.Cells.Item(10, 11).Value = 60000000
$ws.Cells.Item(10, 12).Value = 125
$ws.Cells.Item(10, 14).Value = ""

$ws.Cells.Item(11, 11).Value = 60000168
$ws.Cells.Item(11, 14).Value = 60000169

$ws.Cells.Item(12, 11).Value = 60000000
$ws.Cells.Item(12, 12).Value = 125
$ws.Cells.Item(12, 14).Value = 60000169

$ws.Cells.Item(13, 11).Value = 60000000
$ws.Cells.Item(13, 12).Value = 125
$ws.Cells.Item(13, 14).Value = 60000000
$ws.Cells.Item(13, 15).Value = 126

$ws.Cells.Item(14, 11).Value = 60000168
$ws.Cells.Item(14, 14).Value = ""

$ws.Cells.Item(15, 11).Value = 60000000
$ws.Cells.Item(15, 12).Value = 125
$ws.Cells.Item(15, 14).Value = ""

$ws.Cells.Item(16, 11).Value = 60000168
$ws.Cells.Item(16, 14).Value = 60000169

$ws.Cells.Item(17, 11).Value = 60000000
$ws.Cells.Item(17, 12).Value = 125
$ws.Cells.Item(17, 14).Value = 60000169

$ws.Cells.Item(18, 11).Value = 60000000
$ws.Cells.Item(18, 12).Value = 125
$ws.Cells.Item(18, 14).Value = 60000000
$ws.Cells.Item(18, 15).Value = 126

$ws.Cells.Item(19, 11).Value = 60000168
$ws.Cells.Item(19, 14).Value = ""

$ws.Cells.Item(20, 11).Value = 60000000
$ws.Cells.Item(20, 12).Value = 125
$ws.Cells.Item(20, 14).Value = ""

$ws.Cells.Item(21, 11).Value = 60000168
$ws.Cells.Item(21, 14).Value = 60000169

$ws.Cells.Item(22, 11).Value = 60000000
$ws.Cells.Item(22, 12).Value = 125
$ws.Cells.Item(22, 14).Value = 60000169

$ws.Cells.Item(23, 11).Value = 60000000
$ws.Cells.Item(23, 12).Value = 125
$ws.Cells.Item(23, 14).Value = 60000000
$ws.Cells.Item(23, 15).Value = 126

$ws.Cells.Item(24, 11).Value = 60000168
$ws.Cells.Item(24, 14).Value = ""

$ws.Cells.Item(25, 11).Value = 60000000
$ws.Cells.Item(25, 12).Value = 125
$ws.Cells.Item(25, 14).Value = ""

$ws.Cells.Item(26, 11).Value = 60000168
$ws.Cells.Item(26, 14).Value = 60000169

$ws.Cells.Item(27, 11).Value = 60000000
$ws.Cells.Item(27, 12).Value = 125
$ws.Cells.Item(27, 14).Value = 60000169

$ws.Cells.Item(28, 11).Value = 60000000
$ws.Cells.Item(28, 12).Value = 125
$ws.Cells.Item(28, 14).Value = 60000000
$ws.Cells.Item(28, 15).Value = 126

$ws.Cells.Item(29, 11).Value = 60000168
$ws.Cells.Item(29, 14).Value = ""
